# Changes of 5th May 2022
# Replace the PackageTrackNum values in column C (rows 2-22) with a new
# batch of tracking numbers, also updating column D where it mirrors
# column C (rows 5, 6, 7, 13, 14, 15, 16, 17).
#
# The track numbers are long all-digit strings that must stay text (as in
# the original file), so each value is written with a leading apostrophe
# (forces text entry) and the cell style is then reset back to "Normal" so
# the cell keeps its original (default) formatting instead of picking up
# the quote-prefix style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTrackNums = @{
    2  = "320018471602"
    3  = "320018471613"
    4  = "320018471646"
    5  = "320018471668"
    6  = "320018471705"
    7  = "320018471727"
    8  = "320018471750"
    9  = "320018471771"
    10 = "320018471808"
    11 = "320018471820"
    12 = "320018471863"
    13 = "320018471885"
    14 = "320018471911"
    15 = "320018471933"
    16 = "320018471966"
    17 = "320018471988"
    18 = "320018472024"
    19 = "320018472046"
    20 = "320018472079"
    21 = "320018472090"
    22 = "320018472127"
}

# Rows where column D mirrors column C's track number.
$mirrorDRows = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in 2..22) {
    $value = $newTrackNums[$row]

    $cRange = $ws.Range("C$row")
    $cRange.Value = "'" + $value
    $cRange.Style = "Normal"

    if ($mirrorDRows -contains $row) {
        $dRange = $ws.Range("D$row")
        $dRange.Value = "'" + $value
        $dRange.Style = "Normal"
    }
}
